# Edit script: apply data refresh (2024 dataset), column width tweaks, and remove last row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update numeric data in rows 2-5 (new 1000-row-batch sample values) ---
$ws.Cells.Item(2,1).Value = 45116.50694444445
$ws.Cells.Item(2,2).Value = 23.541
$ws.Cells.Item(2,3).Value = 16.304
$ws.Cells.Item(2,4).Value = 4.244
$ws.Cells.Item(2,5).Value = 49.627
$ws.Cells.Item(2,6).Value = 41.058
$ws.Cells.Item(2,7).Value = 18.526
$ws.Cells.Item(2,8).Value = 61.718
$ws.Cells.Item(2,9).Value = 28.505
$ws.Cells.Item(2,10).Value = 12.158
$ws.Cells.Item(2,11).Value = 18.744
$ws.Cells.Item(2,12).Value = 19.59
$ws.Cells.Item(2,13).Value = 20.45
$ws.Cells.Item(2,14).Value = 5.915
$ws.Cells.Item(2,15).Value = 18.422
$ws.Cells.Item(2,16).Value = 25.954
$ws.Cells.Item(2,17).Value = 15.331
$ws.Cells.Item(2,18).Value = 3.814
$ws.Cells.Item(2,19).Value = 2.604
$ws.Cells.Item(2,20).Value = 273.066
$ws.Cells.Item(2,21).Value = 51.313
$ws.Cells.Item(2,22).Value = 17.004
$ws.Cells.Item(2,23).Value = 34.141
$ws.Cells.Item(2,24).Value = 17.76
$ws.Cells.Item(2,25).Value = 2.295
$ws.Cells.Item(2,26).Value = 30.618
$ws.Cells.Item(2,27).Value = 15.02
$ws.Cells.Item(2,28).Value = 13.452
$ws.Cells.Item(2,29).Value = 15.735
$ws.Cells.Item(2,30).Value = 20.304
$ws.Cells.Item(2,31).Value = 3.64
$ws.Cells.Item(2,32).Value = 54.593
$ws.Cells.Item(2,33).Value = 9.503
$ws.Cells.Item(2,34).Value = 21.259
$ws.Cells.Item(3,1).Value = 45116.51388888889
$ws.Cells.Item(3,2).Value = 1.922
$ws.Cells.Item(3,3).Value = 0.766
$ws.Cells.Item(3,4).Value = 1.33
$ws.Cells.Item(3,5).Value = 3.557
$ws.Cells.Item(3,6).Value = 2.984
$ws.Cells.Item(3,7).Value = 1.517
$ws.Cells.Item(3,8).Value = 13.693
$ws.Cells.Item(3,9).Value = 2.327
$ws.Cells.Item(3,10).Value = 0.881
$ws.Cells.Item(3,11).Value = 1.293
$ws.Cells.Item(3,12).Value = 1.391
$ws.Cells.Item(3,13).Value = 1.294
$ws.Cells.Item(3,14).Value = 0.515
$ws.Cells.Item(3,15).Value = 1.504
$ws.Cells.Item(3,16).Value = 2.133
$ws.Cells.Item(3,17).Value = 1.532
$ws.Cells.Item(3,18).Value = 1.511
$ws.Cells.Item(3,19).Value = 0.642
$ws.Cells.Item(3,20).Value = 15.714
$ws.Cells.Item(3,21).Value = 4.787
$ws.Cells.Item(3,22).Value = 1.388
$ws.Cells.Item(3,23).Value = 3.007
$ws.Cells.Item(3,24).Value = 1.618
$ws.Cells.Item(3,25).Value = 0.095
$ws.Cells.Item(3,26).Value = 5.779
$ws.Cells.Item(3,27).Value = 1.226
$ws.Cells.Item(3,28).Value = 1.343
$ws.Cells.Item(3,29).Value = 1.51
$ws.Cells.Item(3,30).Value = 1.355
$ws.Cells.Item(3,31).Value = 1.294
$ws.Cells.Item(3,32).Value = 12.829
$ws.Cells.Item(3,33).Value = 0.641
$ws.Cells.Item(3,34).Value = 1.746
$ws.Cells.Item(4,1).Value = 45116.52083333334
$ws.Cells.Item(4,2).Value = 22.58
$ws.Cells.Item(4,3).Value = 16.515
$ws.Cells.Item(4,4).Value = 1.53
$ws.Cells.Item(4,5).Value = 48.714
$ws.Cells.Item(4,6).Value = 40.315
$ws.Cells.Item(4,7).Value = 17.77
$ws.Cells.Item(4,8).Value = 64.102
$ws.Cells.Item(4,9).Value = 27.341
$ws.Cells.Item(4,10).Value = 12.06
$ws.Cells.Item(4,11).Value = 18.137
$ws.Cells.Item(4,12).Value = 19.576
$ws.Cells.Item(4,13).Value = 20.549
$ws.Cells.Item(4,14).Value = 5.674
$ws.Cells.Item(4,15).Value = 17.67
$ws.Cells.Item(4,16).Value = 25.084
$ws.Cells.Item(4,17).Value = 14.895
$ws.Cells.Item(4,18).Value = 1.14
$ws.Cells.Item(4,19).Value = 0.995
$ws.Cells.Item(4,20).Value = 261.627
$ws.Cells.Item(4,21).Value = 49.28
$ws.Cells.Item(4,22).Value = 16.31
$ws.Cells.Item(4,23).Value = 33.07
$ws.Cells.Item(4,24).Value = 17.606
$ws.Cells.Item(4,25).Value = 2.228
$ws.Cells.Item(4,26).Value = 31.618
$ws.Cells.Item(4,27).Value = 14.407
$ws.Cells.Item(4,28).Value = 12.825
$ws.Cells.Item(4,29).Value = 15.059
$ws.Cells.Item(4,30).Value = 20.516
$ws.Cells.Item(4,31).Value = 0.784
$ws.Cells.Item(4,32).Value = 57.685
$ws.Cells.Item(4,33).Value = 9.18
$ws.Cells.Item(4,34).Value = 20.391
$ws.Cells.Item(5,1).Value = 45116.52777777778
$ws.Cells.Item(5,2).Value = 2.88
$ws.Cells.Item(5,3).Value = 1.83
$ws.Cells.Item(5,4).Value = 0.65
$ws.Cells.Item(5,5).Value = 5.97
$ws.Cells.Item(5,6).Value = 4.96
$ws.Cells.Item(5,7).Value = 2.27
$ws.Cells.Item(5,8).Value = 16.21
$ws.Cells.Item(5,9).Value = 3.49
$ws.Cells.Item(5,10).Value = 1.52
$ws.Cells.Item(5,11).Value = 2.17
$ws.Cells.Item(5,12).Value = 2.45
$ws.Cells.Item(5,13).Value = 2.43
$ws.Cells.Item(5,14).Value = 0.75
$ws.Cells.Item(5,15).Value = 2.26
$ws.Cells.Item(5,16).Value = 3.25
$ws.Cells.Item(5,17).Value = 2.09
$ws.Cells.Item(5,18).Value = 0.72
$ws.Cells.Item(5,19).Value = 0.31
$ws.Cells.Item(5,20).Value = 27.12
$ws.Cells.Item(5,21).Value = 6.78
$ws.Cells.Item(5,22).Value = 2.08
$ws.Cells.Item(5,23).Value = 4.47
$ws.Cells.Item(5,24).Value = 2.4
$ws.Cells.Item(5,25).Value = 0.22
$ws.Cells.Item(5,26).Value = 7.14
$ws.Cells.Item(5,27).Value = 1.84
$ws.Cells.Item(5,28).Value = 1.77
$ws.Cells.Item(5,29).Value = 2.04
$ws.Cells.Item(5,30).Value = 2.52
$ws.Cells.Item(5,31).Value = 0.5600000000000001
$ws.Cells.Item(5,32).Value = 15.19
$ws.Cells.Item(5,33).Value = 1.08
$ws.Cells.Item(5,34).Value = 2.61

# Update column widths
$ws.Columns(2).ColumnWidth = 7.15
$ws.Columns(3).ColumnWidth = 7.15
$ws.Columns(6).ColumnWidth = 7.15
$ws.Columns(7).ColumnWidth = 7.15
$ws.Columns(9).ColumnWidth = 7.15
$ws.Columns(10).ColumnWidth = 7.15
$ws.Columns(11).ColumnWidth = 7.15
$ws.Columns(12).ColumnWidth = 7.15
$ws.Columns(13).ColumnWidth = 7.15
$ws.Columns(15).ColumnWidth = 7.15
$ws.Columns(16).ColumnWidth = 7.15
$ws.Columns(17).ColumnWidth = 7.15
$ws.Columns(22).ColumnWidth = 7.15
$ws.Columns(23).ColumnWidth = 7.15
$ws.Columns(24).ColumnWidth = 7.15
$ws.Columns(26).ColumnWidth = 7.15
$ws.Columns(27).ColumnWidth = 7.15
$ws.Columns(28).ColumnWidth = 7.15
$ws.Columns(29).ColumnWidth = 7.15
$ws.Columns(30).ColumnWidth = 7.15
$ws.Columns(32).ColumnWidth = 7.15
$ws.Columns(34).ColumnWidth = 7.15
$ws.Columns(20).ColumnWidth = 8.15

# --- Remove row 6 (dataset now only has 4 data rows) ---
$ws.Rows(6).Delete()
